$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---------------------------------------------------------------------------
# Row 1 (headers).
# The first three header cells keep the same meaning as "bank / deposit_type /
# currency" column captions; columns E/F, which used to (incorrectly) carry a
# sample owner/amount in the header row, become the real "owner"/"total"
# headers, and seven more metadata headers are appended, matching the layout
# already used on the other sheets (土地/建物/汽車/...).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Give the new header cells (G1:M1) the same bold / bordered look as the rest
# of row 1.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Data rows 2-5.
# Bank / deposit type / currency stay as before; the new metadata columns
# (property_category .. index) are appended with the same constant values used
# on every other sheet of this workbook.
# The "date" column holds a literal text value ("2011-12-27"), not a real
# Excel date, so force a text number format before assigning it to stop Excel
# from reinterpreting the string as a date serial.
# ---------------------------------------------------------------------------
$ws.Range("I2:I5").NumberFormat = "@"

$ws.Range("B2").Value = "台中商業銀行台中港分行"
$ws.Range("C2").Value = "定期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2011-12-27"
$ws.Range("J2").Value = "楊瓊瓔"
$ws.Range("K2").Value = 854
$ws.Range("L2").Value = "tmpd1401"
$ws.Range("M2").Value = 58

$ws.Range("B3").Value = "臺灣銀行台中港分行"
$ws.Range("C3").Value = "定期儲蓄存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2011-12-27"
$ws.Range("J3").Value = "楊瓊瓔"
$ws.Range("K3").Value = 854
$ws.Range("L3").Value = "tmpd1401"
$ws.Range("M3").Value = 59

$ws.Range("B4").Value = "臺灣銀行台中港分行"
$ws.Range("C4").Value = "綜合存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2011-12-27"
$ws.Range("J4").Value = "楊瓊瓔"
$ws.Range("K4").Value = 854
$ws.Range("L4").Value = "tmpd1401"
$ws.Range("M4").Value = 60

$ws.Range("B5").Value = "台中商業銀行台中港分行"
$ws.Range("C5").Value = "綜合存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2011-12-27"
$ws.Range("J5").Value = "楊瓊瓔"
$ws.Range("K5").Value = 854
$ws.Range("L5").Value = "tmpd1401"
$ws.Range("M5").Value = 61

# Apply the plain (non-bold, no border) data-row style used by column B..F to
# the new metadata columns so G2:M5 visually match the rest of the table.
$ws.Range("B2").Copy()
$ws.Range("G2:M5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
